# Swap the data rows for "Ka'imi Fairbairn" (rows 2-4) and "BrandonMcManus" (rows 11-13).
# Only columns A (Player), C (FG%), D (Lng) change; column B (Season Group) stays put.
# NOTE: reading must use .Value2 (plain .Value getter returns a property descriptor
# in this COM-interop shim); writing works fine via plain .Value = ... assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for both blocks before overwriting anything.
$origTop = @()
for ($r = 2; $r -le 4; $r++) {
    $origTop += , @($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 3).Value2, $ws.Cells.Item($r, 4).Value2)
}

$origBottom = @()
for ($r = 11; $r -le 13; $r++) {
    $origBottom += , @($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 3).Value2, $ws.Cells.Item($r, 4).Value2)
}

# Write bottom block's (BrandonMcManus) values into the top block (rows 2-4).
for ($i = 0; $i -lt 3; $i++) {
    $r = 2 + $i
    $vals = $origBottom[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
}

# Write top block's (Ka'imi Fairbairn) values into the bottom block (rows 11-13).
for ($i = 0; $i -lt 3; $i++) {
    $r = 11 + $i
    $vals = $origTop[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
}
